$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Density"
$ws.Range("B3").Value = 2000
$ws.Range("C3").Value = "kg/m3"

$ws.Range("C3").Select() | Out-Null
